$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-10 Wednesday" "2024-01-11 Thursday"

Replace-Text "574×7=" "307×4="
Replace-Text "459×7=" "402×7="
Replace-Text "385×7=" "344×2="
Replace-Text "134×8=" "618×7="
Replace-Text "445×6=" "239×6="

Replace-Text "198×6=" "185×4="
Replace-Text "237×7=" "732×8="
Replace-Text "822×6=" "647×8="
Replace-Text "599×8=" "196×4="
Replace-Text "806×6=" "667×3="

Replace-Text "740×5=" "406×4="
Replace-Text "649×6=" "229×8="
Replace-Text "719×6=" "424×3="
Replace-Text "303×5=" "541×9="
Replace-Text "296×8=" "511×3="

Replace-Text "319×7=" "305×6="
Replace-Text "696×5=" "230×8="
Replace-Text "728×2=" "964×2="
Replace-Text "548×8=" "508×5="
Replace-Text "343×8=" "624×9="

Replace-Text "604×3=" "693×8="
Replace-Text "788×7=" "844×4="
Replace-Text "894×9=" "304×9="
Replace-Text "356×5=" "895×8="
Replace-Text "527×7=" "597×3="
